$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Number of beers" header (column L) to just "Beers"
$ws.Range("L2").Value = "Beers"

# Shift the existing two data-row dates forward by one day
$ws.Range("A3").Value = 43791
$ws.Range("A4").Value = 43792

# Add a new day's data as row 5
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("A5").Value = 43793
$ws.Range("B5").Value = 4
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 16.5
$ws.Range("G5").Value = 5
$ws.Range("H5").Value = 20
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 3
$ws.Range("M5").Value = 3

# Move the selection to M2, matching the saved workbook's UI state
$null = $ws.Range("M2").Select()
